$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 6962
$ws.Range("I96").Value = 6962
$ws.Range("K96").Value = 20886
$ws.Range("M96").Value = -19513
$ws.Range("H117").Value = 96219.75
$ws.Range("J117").Value = 96219.75
$ws.Range("L117").Value = 96219.75
$ws.Range("N117").Value = -105397.75
$ws.Range("H138").Value = 5290.195
$ws.Range("J138").Value = 5325.127
$ws.Range("L138").Value = 15975.381
$ws.Range("N138").Value = -26255.381
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 754.6667
$ws.Range("I4").Value = 751.875
$ws.Range("K4").Value = 751.875
$ws.Range("M4").Value = -635.875
$ws.Range("H8").Value = 6649.5
$ws.Range("J8").Value = 12299
$ws.Range("L8").Value = 12299
$ws.Range("N8").Value = -12587
$ws.Range("H32").Value = 40033.95
$ws.Range("I32").Value = 35512.176
$ws.Range("J32").Value = 72978.28999999999
$ws.Range("K32").Value = 35512.176
$ws.Range("L32").Value = 72978.28999999999
$ws.Range("M32").Value = -35225.176
$ws.Range("N32").Value = -73552.28999999999
$ws.Range("H102").Value = 40629
$ws.Range("I102").Value = 40629
$ws.Range("K102").Value = 40629
$ws.Range("M102").Value = -39007
$ws.Range("H122").Value = 3198.25
$ws.Range("I122").Value = 2770
$ws.Range("K122").Value = 8310
$ws.Range("M122").Value = -5860
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 93724.5
$ws.Range("J9").Value = 93724.5
$ws.Range("L9").Value = 93724.5
$ws.Range("N9").Value = -94060.5
$ws.Range("H134").Value = 4506625
$ws.Range("I134").Value = 2098.7144
$ws.Range("K134").Value = 6296.1432
$ws.Range("M134").Value = -3761.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 642187.25
$ws.Range("I6").Value = 847936.5
$ws.Range("K6").Value = 847936.5
$ws.Range("M6").Value = -847823.5
$ws.Range("H12").Value = 85
$ws.Range("I12").Value = 85
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 85
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 85
$ws.Range("N12").ClearContents()
$ws.Range("H22").Value = 3500
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5700
$ws.Range("H31").Value = 23837.2
$ws.Range("I31").Value = 51233.332
$ws.Range("J31").Value = 12096
$ws.Range("K31").Value = 51233.332
$ws.Range("L31").Value = 12096
$ws.Range("M31").Value = -50938.332
$ws.Range("N31").Value = -12686
$ws.Range("H34").Value = 23837.2
$ws.Range("I34").Value = 51233.332
$ws.Range("J34").Value = 12096
$ws.Range("K34").Value = 51233.332
$ws.Range("L34").Value = 12096
$ws.Range("M34").Value = -51031.332
$ws.Range("N34").Value = -12500
$ws.Range("H58").Value = 4027.1538
$ws.Range("I58").Value = 1449.8
$ws.Range("J58").Value = 5638
$ws.Range("K58").Value = 1449.8
$ws.Range("L58").Value = 5638
$ws.Range("M58").Value = -1246.8
$ws.Range("N58").Value = -6044
$ws.Range("H136").Value = 4027.1538
$ws.Range("I136").Value = 1449.8
$ws.Range("J136").Value = 5638
$ws.Range("K136").Value = 4349.4
$ws.Range("L136").Value = 16914
$ws.Range("M136").Value = -1799.4
$ws.Range("N136").Value = -22014
$ws.Range("H141").Value = 331045.03
$ws.Range("J141").Value = 338951.94
$ws.Range("L141").Value = 338951.94
$ws.Range("N141").Value = -349311.94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 390861
$ws.Range("J131").Value = 2919.3
$ws.Range("L131").Value = 8757.900000000001
$ws.Range("N131").Value = -18837.9
$ws.Range("H132").Value = 2301.125
$ws.Range("I132").Value = 2168
$ws.Range("K132").Value = 19512
$ws.Range("M132").Value = -16982
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 174.82353
$ws.Range("I2").Value = 71.75
$ws.Range("J2").Value = 266.44446
$ws.Range("K2").Value = 71.75
$ws.Range("L2").Value = 266.44446
$ws.Range("M2").Value = 41.25
$ws.Range("N2").Value = -492.44446
$ws.Range("H13").Value = 19989.5
$ws.Range("I13").Value = 14402
$ws.Range("K13").Value = 14402
$ws.Range("M13").Value = -14263
$ws.Range("H122").Value = 4830.3335
$ws.Range("I122").Value = 4609.25
$ws.Range("J122").Value = 5272.5
$ws.Range("K122").Value = 13827.75
$ws.Range("L122").Value = 15817.5
$ws.Range("M122").Value = -11377.75
$ws.Range("N122").Value = -20717.5
$ws.Range("H132").Value = 5109.9287
$ws.Range("I132").Value = 4093.2727
$ws.Range("J132").Value = 8837.666999999999
$ws.Range("K132").Value = 12279.8181
$ws.Range("L132").Value = 26513.001
$ws.Range("M132").Value = -9749.8181
$ws.Range("N132").Value = -31573.001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6881.68
$ws.Range("I122").Value = 6056.7856
$ws.Range("K122").Value = 18170.3568
$ws.Range("M122").Value = -15720.3568
$ws.Range("H132").Value = 11333.333
$ws.Range("I132").Value = 14500
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 43500
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -40970
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 68240.086
$ws.Range("I136").Value = 20432.076
$ws.Range("J136").Value = 130390.5
$ws.Range("K136").Value = 61296.228
$ws.Range("L136").Value = 391171.5
$ws.Range("M136").Value = -58746.228
$ws.Range("N136").Value = -396271.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28591342
$ws.Range("I62").Value = 7400
$ws.Range("J62").Value = 33355334
$ws.Range("K62").Value = 7400
$ws.Range("L62").Value = 33355334
$ws.Range("M62").Value = -6776
$ws.Range("N62").Value = -33356582
$ws.Range("H65").Value = 28591342
$ws.Range("I65").Value = 7400
$ws.Range("J65").Value = 33355334
$ws.Range("K65").Value = 37000
$ws.Range("L65").Value = 166776670
$ws.Range("M65").Value = -33880
$ws.Range("N65").Value = -166782910
$ws.Range("H81").Value = 2440.7144
$ws.Range("I81").Value = 2440.7144
$ws.Range("K81").Value = 4881.4288
$ws.Range("M81").Value = -3820.4288
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766
$ws.Range("H84").Value = 2440.7144
$ws.Range("I84").Value = 2440.7144
$ws.Range("K84").Value = 24407.144
$ws.Range("M84").Value = -19103.144
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652
$ws.Range("I96").Value = 2363.5
$ws.Range("J96").Value = 2802
$ws.Range("K96").Value = 2363.5
$ws.Range("L96").Value = 2802
$ws.Range("M96").Value = -990.5
$ws.Range("N96").Value = -5548
$ws.Range("H100").Value = 1491.6333
$ws.Range("I100").Value = 1405.9231
$ws.Range("J100").Value = 2048.75
$ws.Range("K100").Value = 2811.8462
$ws.Range("L100").Value = 4097.5
$ws.Range("M100").Value = -2270.8462
$ws.Range("N100").Value = -5179.5
$ws.Range("H132").Value = 3486014
$ws.Range("J132").Value = 8356712
$ws.Range("L132").Value = 25070136
$ws.Range("N132").Value = -25075196
